$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 2 ("Tue 24 Oct"),
# shifting the existing data rows down.
$ws.Range("A2:E5").Insert()

# Fill in the newly inserted rows with the additional days of data.
$ws.Range("A2").Value = "Fri 20 Oct"
$ws.Range("B2").Value = 3257
$ws.Range("C2").Value = 13531
$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 10.73

$ws.Range("A3").Value = "Sat 21 Oct"
$ws.Range("B3").Value = 3167
$ws.Range("C3").Value = 15717
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 11.79

$ws.Range("A4").Value = "Sun 22 Oct"
$ws.Range("B4").Value = 3224
$ws.Range("C4").Value = 14589
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 10.73

$ws.Range("A5").Value = "Mon 23 Oct"
$ws.Range("B5").Value = 2757
$ws.Range("C5").Value = 10216
$ws.Range("D5").Value = 18
$ws.Range("E5").Value = 8.12

# Update today's calories burned figure (row 9, was row 5 before the insert).
$ws.Range("B9").Value = 588
